# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with freshly scraped values. Every target cell already holds
# plain text (t="inlineStr" in the OOXML), e.g. "590.75" or "  +1.45%  ".
#
# Excel's Range.Value setter auto-coerces a numeric-looking string (like
# "590.41") into a real number, which would both change the stored cell
# type and round the text (590.41 -> 590.40999999999997 in floating
# point). To keep these values as literal text -- matching the source
# data exactly -- cells whose new value parses as a plain number are
# first switched to the "@" (Text) number format, written, then
# restored to the "Normal" cell style so no stray style/format diff is
# left on the cell. Cells whose new value is not numeric (extra "."
# separators, a trailing "%" sign, padding spaces, subscript digits,
# etc.) are never auto-converted by Excel, so they are assigned
# directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.447.07"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "2.644.78"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "590.41"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.41%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "143.53"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "2.643.23"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  +0.30%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.33"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "3.120.11"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "63.348.34"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "2.671.81"
$ws.Range("E18").Value = "  +3.48%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.39"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.51%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "339.42"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.30%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.36"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.73"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("E23").Value = "  +0.03%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "67.16"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  +5.35%  "
$ws.Range("E26").Value = "  -0.31%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.52"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.00%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "543.06"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +17.88%  "
$ws.Range("E29").Value = "  +0.34%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.41"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.46%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.75"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("E32").Value = "  +12.32%  "
$ws.Range("E33").Value = "  +2.66%  "
$ws.Range("D34").Value = "0.0₃0805"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "175.24"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  +9.01%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  +0.39%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "19.04"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("E40").Value = "  +6.18%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "172.08"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +9.36%  "
$ws.Range("E42").Value = "  -0.04%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "40.36"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("E44").Value = "  +0.16%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "22.24"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +5.43%  "
$ws.Range("E46").Value = "  +4.96%  "
$ws.Range("E47").Value = "  +0.79%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0960"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("E49").Value = "  +1.78%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "18.70"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("E51").Value = "  -0.58%  "
